$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = 6021; E = 0.04284750337381917 }
    @{ Row = 3; D = 715; E = -0.1452380952380952 }
    @{ Row = 4; D = 3766; E = -0.1130063965884861 }
    @{ Row = 5; D = 1409; E = -0.1034912718204489 }
    @{ Row = 6; D = 3660; E = -0.1753336401288541 }
    @{ Row = 7; D = 3358; E = -0.1197668256491786 }
    @{ Row = 8; D = 3207; E = -0.08436992969172526 }
    @{ Row = 9; D = 3827; E = -0.04039408866995074 }
    @{ Row = 10; D = 1611; E = 0.03775510204081633 }
    @{ Row = 11; D = 2574; E = -0.1005398110661269 }
    @{ Row = 12; D = 1824; E = -0.2561576354679803 }
    @{ Row = 13; D = 873; E = -0.2642706131078224 }
    @{ Row = 14; D = 1444; E = 0.05528255528255528 }
    @{ Row = 15; D = 7176; E = -0.1748421563865954 }
    @{ Row = 16; D = 2519; E = 0.03440702781844802 }
    @{ Row = 17; D = 3845; E = -0.2884895580378825 }
    @{ Row = 18; D = 964; E = -0.06557377049180328 }
    @{ Row = 19; D = 4550; E = -0.1861042183622829 }
    @{ Row = 20; D = 1161; E = -0.2021660649819494 }
    @{ Row = 21; D = 2229; E = -0.0370919881305638 }
    @{ Row = 22; D = 2592; E = -0.07185234014502308 }
    @{ Row = 23; D = 1202; E = 0.02303030303030303 }
    @{ Row = 24; D = 1911; E = -0.1180904522613065 }
    @{ Row = 25; D = 1396; E = 0.02680965147453083 }
    @{ Row = 26; D = 2872; E = -0.1934156378600823 }
    @{ Row = 27; D = 1208; E = -0.07919463087248323 }
    @{ Row = 28; D = 1580; E = -0.2807453416149068 }
    @{ Row = 29; D = 2415; E = -0.3592622293504411 }
    @{ Row = 30; D = 6469; E = 0.04612412460758271 }
    @{ Row = 31; D = 760; E = -0.2048192771084337 }
    @{ Row = 32; D = 2997; E = -0.04973183812774257 }
    @{ Row = 33; D = 1651; E = -0.1559074299634592 }
    @{ Row = 34; D = 2865; E = -0.2178871548619448 }
    @{ Row = 35; D = 1921; E = -0.1150278293135436 }
    @{ Row = 36; D = 1641; E = 0.04869857262804366 }
    @{ Row = 37; D = 2829; E = 0.03104786545924968 }
    @{ Row = 38; D = 726; E = -0.2156862745098039 }
    @{ Row = 39; D = 3266; E = -0.2807799442896936 }
    @{ Row = 40; D = 3367; E = -0.1776251226692836 }
    @{ Row = 41; D = 2083; E = -0.05617977528089887 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 4).Value = $u.D
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
